# "Device" -> "Apparatus" rename across the workbook (IEEE_57Bus example)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Device")

# Rename the sheet itself
$ws.Name = "Apparatus"

# Update the header/label cells (set B2/C2 first so the new shared
# strings land in the same order the author's Excel session produced)
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Make the renamed sheet the active tab/selection, like the author did
# when finishing the edit in Excel
$ws.Activate()
$ws.Range("A2").Select()
